$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (P2, Q2, R2, S2) and the label text in U2
$ws.Range("P2").Value = 0.01
$ws.Range("Q2").Value = -2.432000621081636
$ws.Range("R2").Value = 3.591031790153992
$ws.Range("S2").Value = -8.55457746663083
$ws.Range("U2").Value = "blink+templerun"

# Update row 3 values (Q3, R3, S3) and the label text in U3
$ws.Range("Q3").Value = -2.338056636100122
$ws.Range("R3").Value = 2.220386220614815
$ws.Range("S3").Value = -10.31255766339912
$ws.Range("U3").Value = "blink+sudoku"

# Remove rows 4-6 entirely (data for subjects no longer present)
$ws.Rows("4:6").Delete()
